# Insert a new data row right above current row 383 ("Región de La Araucanía",
# 2022-02-24 record). This shifts the existing row 383 (and everything below
# it) down by one row. The new row 383 is a duplicate of the original row 383
# record, but dated 2023-02-27 (serial 44984) instead of 2022-02-24 (serial
# 44616).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 383 (and below) down by one to make room for the new record.
$ws.Rows(383).Insert()

# Populate the newly inserted row 383 with the duplicated record.
$ws.Cells.Item(383, 1).Value = 10
$ws.Cells.Item(383, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(383, 3).Value = 'La Araucanía'
$ws.Cells.Item(383, 4).Value = 44984
$ws.Cells.Item(383, 5).Value = 9
$ws.Cells.Item(383, 6).Value = 100114013
$ws.Cells.Item(383, 7).Value = 'Zanahoria'
$ws.Cells.Item(383, 8).Value = 'Sin especificar'
$ws.Cells.Item(383, 9).Value = 'Primera'
$ws.Cells.Item(383, 10).Value = 200
$ws.Cells.Item(383, 11).Value = 8000
$ws.Cells.Item(383, 12).Value = 8000
$ws.Cells.Item(383, 13).Value = 8000
$ws.Cells.Item(383, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(383, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(383, 16).Value = 320
$ws.Cells.Item(383, 17).Value = 25
$ws.Cells.Item(383, 18).Value = 'Hortaliza'
